$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D ("Price") hold numeric-looking text (e.g. "1.00",
# "0.600", "76.393.35"). Excel would silently reinterpret a plain
# assignment of these strings as numbers (dropping formatting such as
# trailing zeros), so each Price cell we touch is switched to the Text
# number format first to preserve the literal text exactly.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '76.393.35'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.078.11'
$ws.Range('E3').Value = '  +5.18%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '198.86'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '618.10'
$ws.Range('E6').Value = '  +3.99%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.214'
$ws.Range('E8').Value = '  +8.53%  '
$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.552'
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '3.079.04'
$ws.Range('E10').Value = '  +5.24%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.447'
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.26'
$ws.Range('E13').Value = '  +7.71%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.643.35'
$ws.Range('E14').Value = '  +5.22%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '29.37'
$ws.Range('E15').Value = '  +4.15%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000197'
$ws.Range('E16').Value = '  +4.42%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '76.234.96'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.086.58'
$ws.Range('E18').Value = '  +5.69%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.52'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '9.03'
$ws.Range('E20').Value = '  +3.48%  '
$ws.Range('E21').Value = '  +15.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '384.10'
$ws.Range('E22').Value = '  +2.90%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.52'
$ws.Range('E23').Value = '  +5.33%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.44'
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.58'
$ws.Range('E25').Value = '  +8.02%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.238.41'
$ws.Range('E26').Value = '  +5.24%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '72.49'
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.10'
$ws.Range('E29').Value = '  +5.06%  '
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.32'
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('E33').Value = '  +4.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '501.63'
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('E35').Value = '  +5.87%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.128'
$ws.Range('E36').Value = '  +15.44%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '20.92'
$ws.Range('E38').Value = '  +4.13%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '163.31'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '194.72'
$ws.Range('E40').Value = '  +9.22%  '
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('E43').Value = '  -6.55%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('E45').Value = '  +22.37%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.17'
$ws.Range('E46').Value = '  +5.65%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.26'
$ws.Range('E47').Value = '  +7.19%  '
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.46'
$ws.Range('E49').Value = '  +6.26%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '40.85'
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.600'
$ws.Range('E51').Value = '  +1.86%  '
